$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 44, shifting existing rows 44-50 down to 45-51
$ws.Rows("44:44").Insert()

# Populate the newly inserted row 44 with the new weekly data point
$ws.Range("A44").Value = 7
$ws.Range("B44").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C44").Value = "Ñuble"
$ws.Range("D44").Value = 44816
$ws.Range("D44").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E44").Value = 16
$ws.Range("F44").Value = 100112026
$ws.Range("G44").Value = "Haba"
$ws.Range("H44").Value = "Sin especificar"
$ws.Range("I44").Value = "Primera"
$ws.Range("J44").Value = 80
$ws.Range("K44").Value = 12000
$ws.Range("L44").Value = 13000
$ws.Range("M44").Value = 12500
$ws.Range("N44").Value = "$/saco 25 kilos"
$ws.Range("O44").Value = "Provincia del Elquí"
$ws.Range("P44").Value = 500
$ws.Range("Q44").Value = 25
$ws.Range("R44").Value = "Hortaliza"
